$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("three_line")
$ws.Cells.Item(59, 1).Value = "AVADHSUGAR.NS"
$ws.Cells.Item(59, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(59, 3).Value = "hour"
$ws.Cells.Item(59, 4).Value = 45411.42708333334
$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 5).Value = 620
$ws.Cells.Item(59, 6).Value = 45425.63541666666
$ws.Cells.Item(59, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 7).Value = 589.0499877929688
$ws.Cells.Item(59, 8).Value = 45433.38541666666
$ws.Cells.Item(59, 8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 9).Value = 579.9500122070312
$ws.Cells.Item(59, 10).Value = "High"
$ws.Cells.Item(59, 11).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(60, 1).Value = "UGARSUGAR.NS"
$ws.Cells.Item(60, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(60, 3).Value = "hour"
$ws.Cells.Item(60, 4).Value = 45404.38541666666
$ws.Cells.Item(60, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 5).Value = 83.44999694824219
$ws.Cells.Item(60, 6).Value = 45426.38541666666
$ws.Cells.Item(60, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 7).Value = 77
$ws.Cells.Item(60, 8).Value = 45427.38541666666
$ws.Cells.Item(60, 8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 9).Value = 77
$ws.Cells.Item(60, 10).Value = "High"
$ws.Cells.Item(60, 11).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(61, 1).Value = "MOL.NS"
$ws.Cells.Item(61, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(61, 3).Value = "hour"
$ws.Cells.Item(61, 4).Value = 45427.55208333334
$ws.Cells.Item(61, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(61, 5).Value = 85.34999847412109
$ws.Cells.Item(61, 6).Value = 45436.51041666666
$ws.Cells.Item(61, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(61, 7).Value = 82.34999847412109
$ws.Cells.Item(61, 8).Value = 45442.38541666666
$ws.Cells.Item(61, 8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(61, 9).Value = 80.40000152587891
$ws.Cells.Item(61, 10).Value = "High"
$ws.Cells.Item(61, 11).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(62, 1).Value = "GOODYEAR.BO"
$ws.Cells.Item(62, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(62, 3).Value = "hour"
$ws.Cells.Item(62, 4).Value = 45411.42708333334
$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 5).Value = 1160
$ws.Cells.Item(62, 6).Value = 45419.38541666666
$ws.Cells.Item(62, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 7).Value = 1157.099975585938
$ws.Cells.Item(62, 8).Value = 45419.51041666666
$ws.Cells.Item(62, 8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 9).Value = 1157.099975585938
$ws.Cells.Item(62, 10).Value = "Low"
$ws.Cells.Item(62, 11).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(63, 1).Value = "WEBELSOLAR.NS"
$ws.Cells.Item(63, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(63, 3).Value = "hour"
$ws.Cells.Item(63, 4).Value = 45421.38541666666
$ws.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 5).Value = 684.7999877929688
$ws.Cells.Item(63, 6).Value = 45427.42708333334
$ws.Cells.Item(63, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 7).Value = 659.8499755859375
$ws.Cells.Item(63, 8).Value = 45427.51041666666
$ws.Cells.Item(63, 8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 9).Value = 659.8499755859375
$ws.Cells.Item(63, 10).Value = "High"
$ws.Cells.Item(63, 11).Value = "03/06/2024 10:42:03"

$ws = $wb.Worksheets.Item("two_line")
$ws.Cells.Item(33, 1).Value = "RAJSREESUG.NS"
$ws.Cells.Item(33, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(33, 3).Value = "hour"
$ws.Cells.Item(33, 4).Value = 45428.63541666666
$ws.Cells.Item(33, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 5).Value = 61.29999923706055
$ws.Cells.Item(33, 6).Value = 45434.38541666666
$ws.Cells.Item(33, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 7).Value = 61.40000152587891
$ws.Cells.Item(33, 8).Value = "Low"
$ws.Cells.Item(33, 9).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(34, 1).Value = "CAMLINFINE.NS"
$ws.Cells.Item(34, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(34, 3).Value = "hour"
$ws.Cells.Item(34, 4).Value = 45435.63541666666
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 93.65000152587891
$ws.Cells.Item(34, 6).Value = 45443.38541666666
$ws.Cells.Item(34, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 7).Value = 93.80000305175781
$ws.Cells.Item(34, 8).Value = "Low"
$ws.Cells.Item(34, 9).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(35, 1).Value = "CAMLINFINE.NS"
$ws.Cells.Item(35, 2).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(35, 3).Value = "hour"
$ws.Cells.Item(35, 4).Value = 45443.38541666666
$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 93.80000305175781
$ws.Cells.Item(35, 6).Value = 45443.42708333334
$ws.Cells.Item(35, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 7).Value = 93.80000305175781
$ws.Cells.Item(35, 8).Value = "Low"
$ws.Cells.Item(35, 9).Value = "03/06/2024 10:42:03"

$ws = $wb.Worksheets.Item("ph_pl_breakout_line")
$ws.Cells.Item(282, 1).Value = "UPL.NS"
$ws.Cells.Item(282, 2).Value = 45439.59375
$ws.Cells.Item(282, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(282, 3).Value = 529.6500244140625
$ws.Cells.Item(282, 4).Value = 524
$ws.Cells.Item(282, 5).Value = 525.0499877929688
$ws.Cells.Item(282, 6).Value = "High"
$ws.Cells.Item(282, 7).Value = 529.6500244140625
$ws.Cells.Item(282, 8).Value = "hour"
$ws.Cells.Item(282, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(282, 10).Value = 530.5
$ws.Cells.Item(282, 11).Value = 529
$ws.Cells.Item(282, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(283, 1).Value = "KOTAKBANK.NS"
$ws.Cells.Item(283, 2).Value = 45439.51041666666
$ws.Cells.Item(283, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(283, 3).Value = 1722
$ws.Cells.Item(283, 4).Value = 1717
$ws.Cells.Item(283, 5).Value = 1718.349975585938
$ws.Cells.Item(283, 6).Value = "High"
$ws.Cells.Item(283, 7).Value = 1722
$ws.Cells.Item(283, 8).Value = "hour"
$ws.Cells.Item(283, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(283, 10).Value = 1725
$ws.Cells.Item(283, 11).Value = 1719.150024414062
$ws.Cells.Item(283, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(284, 1).Value = "FEDERALBNK.NS"
$ws.Cells.Item(284, 2).Value = 45433.51041666666
$ws.Cells.Item(284, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(284, 3).Value = 165.1999969482422
$ws.Cells.Item(284, 4).Value = 164.1000061035156
$ws.Cells.Item(284, 5).Value = 164.8500061035156
$ws.Cells.Item(284, 6).Value = "High"
$ws.Cells.Item(284, 7).Value = 165.1999969482422
$ws.Cells.Item(284, 8).Value = "hour"
$ws.Cells.Item(284, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(284, 10).Value = 165.5
$ws.Cells.Item(284, 11).Value = 165.1499938964844
$ws.Cells.Item(284, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(285, 1).Value = "FEDERALBNK.NS"
$ws.Cells.Item(285, 2).Value = 45439.55208333334
$ws.Cells.Item(285, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(285, 3).Value = 165.1999969482422
$ws.Cells.Item(285, 4).Value = 164.1499938964844
$ws.Cells.Item(285, 5).Value = 164.6000061035156
$ws.Cells.Item(285, 6).Value = "High"
$ws.Cells.Item(285, 7).Value = 165.1999969482422
$ws.Cells.Item(285, 8).Value = "hour"
$ws.Cells.Item(285, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(285, 10).Value = 165.5
$ws.Cells.Item(285, 11).Value = 165.1499938964844
$ws.Cells.Item(285, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(286, 1).Value = "DALMIASUG.NS"
$ws.Cells.Item(286, 2).Value = 45436.55208333334
$ws.Cells.Item(286, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(286, 3).Value = 382.6000061035156
$ws.Cells.Item(286, 4).Value = 377.2999877929688
$ws.Cells.Item(286, 5).Value = 378.5499877929688
$ws.Cells.Item(286, 6).Value = "High"
$ws.Cells.Item(286, 7).Value = 382.6000061035156
$ws.Cells.Item(286, 8).Value = "hour"
$ws.Cells.Item(286, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(286, 10).Value = 383
$ws.Cells.Item(286, 11).Value = 381.9500122070312
$ws.Cells.Item(286, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(287, 1).Value = "MHRIL.NS"
$ws.Cells.Item(287, 2).Value = 45428.55208333334
$ws.Cells.Item(287, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(287, 3).Value = 402.4500122070312
$ws.Cells.Item(287, 4).Value = 398.7999877929688
$ws.Cells.Item(287, 5).Value = 401.9500122070312
$ws.Cells.Item(287, 6).Value = "Low"
$ws.Cells.Item(287, 7).Value = 398.7999877929688
$ws.Cells.Item(287, 8).Value = "hour"
$ws.Cells.Item(287, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(287, 10).Value = 398.2999877929688
$ws.Cells.Item(287, 11).Value = 399
$ws.Cells.Item(287, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(288, 1).Value = "INDHOTEL.NS"
$ws.Cells.Item(288, 2).Value = 45435.38541666666
$ws.Cells.Item(288, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(288, 3).Value = 576.5999755859375
$ws.Cells.Item(288, 4).Value = 568.1500244140625
$ws.Cells.Item(288, 5).Value = 569.9500122070312
$ws.Cells.Item(288, 6).Value = "High"
$ws.Cells.Item(288, 7).Value = 576.5999755859375
$ws.Cells.Item(288, 8).Value = "hour"
$ws.Cells.Item(288, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(288, 10).Value = 582.5
$ws.Cells.Item(288, 11).Value = 575.5499877929688
$ws.Cells.Item(288, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(289, 1).Value = "AETHER.NS"
$ws.Cells.Item(289, 2).Value = 45426.42708333334
$ws.Cells.Item(289, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(289, 3).Value = 820.5
$ws.Cells.Item(289, 4).Value = 808.8499755859375
$ws.Cells.Item(289, 5).Value = 817.75
$ws.Cells.Item(289, 6).Value = "Low"
$ws.Cells.Item(289, 7).Value = 808.8499755859375
$ws.Cells.Item(289, 8).Value = "hour"
$ws.Cells.Item(289, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(289, 10).Value = 808
$ws.Cells.Item(289, 11).Value = 813.75
$ws.Cells.Item(289, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(290, 1).Value = "RAJSREESUG.NS"
$ws.Cells.Item(290, 2).Value = 45428.63541666666
$ws.Cells.Item(290, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(290, 3).Value = 62
$ws.Cells.Item(290, 4).Value = 61.29999923706055
$ws.Cells.Item(290, 5).Value = 61.54999923706055
$ws.Cells.Item(290, 6).Value = "Low"
$ws.Cells.Item(290, 7).Value = 61.29999923706055
$ws.Cells.Item(290, 8).Value = "hour"
$ws.Cells.Item(290, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(290, 10).Value = 61.09999847412109
$ws.Cells.Item(290, 11).Value = 61.79999923706055
$ws.Cells.Item(290, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(291, 1).Value = "RAJSREESUG.NS"
$ws.Cells.Item(291, 2).Value = 45434.38541666666
$ws.Cells.Item(291, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(291, 3).Value = 63.5
$ws.Cells.Item(291, 4).Value = 61.40000152587891
$ws.Cells.Item(291, 5).Value = 62.25
$ws.Cells.Item(291, 6).Value = "Low"
$ws.Cells.Item(291, 7).Value = 61.40000152587891
$ws.Cells.Item(291, 8).Value = "hour"
$ws.Cells.Item(291, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(291, 10).Value = 61.09999847412109
$ws.Cells.Item(291, 11).Value = 61.79999923706055
$ws.Cells.Item(291, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(292, 1).Value = "CAMLINFINE.NS"
$ws.Cells.Item(292, 2).Value = 45443.38541666666
$ws.Cells.Item(292, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(292, 3).Value = 94.69999694824219
$ws.Cells.Item(292, 4).Value = 93.80000305175781
$ws.Cells.Item(292, 5).Value = 94.30000305175781
$ws.Cells.Item(292, 6).Value = "Low"
$ws.Cells.Item(292, 7).Value = 93.80000305175781
$ws.Cells.Item(292, 8).Value = "hour"
$ws.Cells.Item(292, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(292, 10).Value = 93.25
$ws.Cells.Item(292, 11).Value = 93.84999847412109
$ws.Cells.Item(292, 12).Value = "03/06/2024 10:42:03"
$ws.Cells.Item(293, 1).Value = "CAMLINFINE.NS"
$ws.Cells.Item(293, 2).Value = 45443.42708333334
$ws.Cells.Item(293, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(293, 3).Value = 94.65000152587891
$ws.Cells.Item(293, 4).Value = 93.80000305175781
$ws.Cells.Item(293, 5).Value = 94.30000305175781
$ws.Cells.Item(293, 6).Value = "Low"
$ws.Cells.Item(293, 7).Value = 93.80000305175781
$ws.Cells.Item(293, 8).Value = "hour"
$ws.Cells.Item(293, 9).Value = "03-06-2024 15:15:00"
$ws.Cells.Item(293, 10).Value = 93.25
$ws.Cells.Item(293, 11).Value = 93.84999847412109
$ws.Cells.Item(293, 12).Value = "03/06/2024 10:42:03"

